# Auto-generated Excel COM-interop script
# Applies cell-value updates to match the target diff for Sheets/Unicorn_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 76
$ws_ALC.Range("H76").Value2 = 8195.75
$ws_ALC.Range("I76").Value2 = 11062.929
$ws_ALC.Range("J76").Value2 = 5328.5713
$ws_ALC.Range("K76").Value2 = 11062.929
$ws_ALC.Range("L76").Value2 = 5328.5713
$ws_ALC.Range("M76").Value2 = -10747.929
$ws_ALC.Range("N76").Value2 = -5958.5713

# ALC row 79
$ws_ALC.Range("H79").Value2 = 8195.75
$ws_ALC.Range("I79").Value2 = 11062.929
$ws_ALC.Range("J79").Value2 = 5328.5713
$ws_ALC.Range("K79").Value2 = 11062.929
$ws_ALC.Range("L79").Value2 = 5328.5713
$ws_ALC.Range("M79").Value2 = -9970.929
$ws_ALC.Range("N79").Value2 = -7512.5713

# ALC row 80
$ws_ALC.Range("H80").Value2 = 23608.857
$ws_ALC.Range("I80").Value2 = 67937.60000000001
$ws_ALC.Range("J80").Value2 = 9756.125
$ws_ALC.Range("K80").Value2 = 203812.8
$ws_ALC.Range("L80").Value2 = 29268.375
$ws_ALC.Range("M80").Value2 = -202814.8
$ws_ALC.Range("N80").Value2 = -31264.375

# ALC row 83
$ws_ALC.Range("H83").Value2 = 23608.857
$ws_ALC.Range("I83").Value2 = 67937.60000000001
$ws_ALC.Range("J83").Value2 = 9756.125
$ws_ALC.Range("K83").Value2 = 611438.4
$ws_ALC.Range("L83").Value2 = 87805.125
$ws_ALC.Range("M83").Value2 = -606446.4
$ws_ALC.Range("N83").Value2 = -97789.125

# ALC row 86
$ws_ALC.Range("H86").Value2 = 3766.7742
$ws_ALC.Range("I86").Value2 = 2707.7273
$ws_ALC.Range("J86").Value2 = 4349.25
$ws_ALC.Range("K86").Value2 = 2707.7273
$ws_ALC.Range("L86").Value2 = 4349.25
$ws_ALC.Range("M86").Value2 = -1584.7273
$ws_ALC.Range("N86").Value2 = -6595.25

# ALC row 89
$ws_ALC.Range("H89").Value2 = 3766.7742
$ws_ALC.Range("I89").Value2 = 2707.7273
$ws_ALC.Range("J89").Value2 = 4349.25
$ws_ALC.Range("K89").Value2 = 13538.6365
$ws_ALC.Range("L89").Value2 = 21746.25
$ws_ALC.Range("M89").Value2 = -7922.636500000001
$ws_ALC.Range("N89").Value2 = -32978.25

# ALC row 106
$ws_ALC.Range("H106").Value2 = 5239
$ws_ALC.Range("I106").Value2 = 4033.6
$ws_ALC.Range("J106").Value2 = 6100
$ws_ALC.Range("K106").Value2 = 4033.6
$ws_ALC.Range("L106").Value2 = 6100
$ws_ALC.Range("M106").Value2 = -3402.6
$ws_ALC.Range("N106").Value2 = -7362

# ALC row 125
$ws_ALC.Range("H125").Value2 = 3008
$ws_ALC.Range("I125").Value2 = 4980
$ws_ALC.Range("J125").Value2 = 1036
$ws_ALC.Range("K125").Value2 = 44820
$ws_ALC.Range("L125").Value2 = 9324
$ws_ALC.Range("M125").Value2 = -42360
$ws_ALC.Range("N125").Value2 = -14244

# ALC row 132
$ws_ALC.Range("H132").Value2 = 4771.081
$ws_ALC.Range("I132").Value2 = 2504.4827
$ws_ALC.Range("K132").Value2 = 7513.4481
$ws_ALC.Range("M132").Value2 = -4983.4481

# ARM row 63
$ws_ARM.Range("H63").Value2 = 6787.6113
$ws_ARM.Range("I63").Value2 = 7698.0835
$ws_ARM.Range("J63").Value2 = 4966.6665
$ws_ARM.Range("K63").Value2 = 7698.0835
$ws_ARM.Range("L63").Value2 = 4966.6665
$ws_ARM.Range("M63").Value2 = -7012.0835
$ws_ARM.Range("N63").Value2 = -6338.6665

# ARM row 66
$ws_ARM.Range("H66").Value2 = 6787.6113
$ws_ARM.Range("I66").Value2 = 7698.0835
$ws_ARM.Range("J66").Value2 = 4966.6665
$ws_ARM.Range("K66").Value2 = 38490.4175
$ws_ARM.Range("L66").Value2 = 24833.3325
$ws_ARM.Range("M66").Value2 = -35058.4175
$ws_ARM.Range("N66").Value2 = -31697.3325

# CRP row 99
$ws_CRP.Range("H99").Value2 = 74558
$ws_CRP.Range("I99").Value2 = 113968
$ws_CRP.Range("J99").Value2 = 3620
$ws_CRP.Range("K99").Value2 = 113968
$ws_CRP.Range("L99").Value2 = 3620
$ws_CRP.Range("M99").Value2 = -112470
$ws_CRP.Range("N99").Value2 = -6616

# CRP row 126
$ws_CRP.Range("H126").Value2 = 74558
$ws_CRP.Range("I126").Value2 = 113968
$ws_CRP.Range("J126").Value2 = 3620
$ws_CRP.Range("K126").Value2 = 341904
$ws_CRP.Range("L126").Value2 = 10860
$ws_CRP.Range("M126").Value2 = -339434
$ws_CRP.Range("N126").Value2 = -15800

# CRP row 134
$ws_CRP.Range("H134").Value2 = 2002.258
$ws_CRP.Range("I134").Value2 = 1379.375
$ws_CRP.Range("J134").Value2 = 2666.6667
$ws_CRP.Range("K134").Value2 = 4138.125
$ws_CRP.Range("L134").Value2 = 8000.000100000001
$ws_CRP.Range("M134").Value2 = -1603.125
$ws_CRP.Range("N134").Value2 = -13070.0001

# CUL row 17
$ws_CUL.Range("H17").Value2 = 957.1429000000001
$ws_CUL.Range("I17").Value2 = 200
$ws_CUL.Range("J17").Value2 = 1525
$ws_CUL.Range("K17").Value2 = 600
$ws_CUL.Range("L17").Value2 = 4575
$ws_CUL.Range("M17").Value2 = -431
$ws_CUL.Range("N17").Value2 = -4913

# CUL row 113
$ws_CUL.Range("H113").Value2 = 647.4286
$ws_CUL.Range("J113").Value2 = 647.4286
$ws_CUL.Range("L113").Value2 = 1942.2858
$ws_CUL.Range("N113").Value2 = -6282.2858

# CUL row 122
$ws_CUL.Range("H122").Value2 = 537.7917
$ws_CUL.Range("I122").Value2 = 369
$ws_CUL.Range("J122").Value2 = 1381.75
$ws_CUL.Range("K122").Value2 = 3321
$ws_CUL.Range("L122").Value2 = 12435.75
$ws_CUL.Range("M122").Value2 = -871
$ws_CUL.Range("N122").Value2 = -17335.75

# CUL row 131
$ws_CUL.Range("H131").Value2 = 2065.2856
$ws_CUL.Range("I131").Value2 = 0
$ws_CUL.Range("J131").Value2 = 2065.2856
$ws_CUL.Range("K131").Value2 = 0
$ws_CUL.Range("L131").Value2 = 6195.8568
$ws_CUL.Range("M131").ClearContents()
$ws_CUL.Range("N131").Value2 = -16275.8568

# GSM row 64
$ws_GSM.Range("H64").Value2 = 271271
$ws_GSM.Range("J64").Value2 = 271271
$ws_GSM.Range("L64").Value2 = 271271
$ws_GSM.Range("N64").Value2 = -271767

# GSM row 67
$ws_GSM.Range("H67").Value2 = 271271
$ws_GSM.Range("J67").Value2 = 271271
$ws_GSM.Range("L67").Value2 = 271271
$ws_GSM.Range("N67").Value2 = -272987

# GSM row 80
$ws_GSM.Range("H80").Value2 = 4704.7095
$ws_GSM.Range("I80").Value2 = 5148
$ws_GSM.Range("J80").Value2 = 2857.6667
$ws_GSM.Range("K80").Value2 = 5148
$ws_GSM.Range("L80").Value2 = 2857.6667
$ws_GSM.Range("M80").Value2 = -4150
$ws_GSM.Range("N80").Value2 = -4853.6667

# GSM row 83
$ws_GSM.Range("H83").Value2 = 4704.7095
$ws_GSM.Range("I83").Value2 = 5148
$ws_GSM.Range("J83").Value2 = 2857.6667
$ws_GSM.Range("K83").Value2 = 25740
$ws_GSM.Range("L83").Value2 = 14288.3335
$ws_GSM.Range("M83").Value2 = -20748
$ws_GSM.Range("N83").Value2 = -24272.3335

# GSM row 102
$ws_GSM.Range("H102").Value2 = 4651.143
$ws_GSM.Range("I102").Value2 = 4846
$ws_GSM.Range("J102").Value2 = 2800
$ws_GSM.Range("K102").Value2 = 4846
$ws_GSM.Range("L102").Value2 = 2800
$ws_GSM.Range("M102").Value2 = -3224
$ws_GSM.Range("N102").Value2 = -6044

# LTW row 7
$ws_LTW.Range("H7").Value2 = 2564.9092
$ws_LTW.Range("I7").Value2 = 2629.8333
$ws_LTW.Range("J7").Value2 = 2487
$ws_LTW.Range("K7").Value2 = 2629.8333
$ws_LTW.Range("L7").Value2 = 2487
$ws_LTW.Range("M7").Value2 = -2517.8333
$ws_LTW.Range("N7").Value2 = -2711

# LTW row 11
$ws_LTW.Range("H11").Value2 = 5200
$ws_LTW.Range("J11").Value2 = 5200
$ws_LTW.Range("L11").Value2 = 5200
$ws_LTW.Range("N11").Value2 = -5480

# LTW row 40
$ws_LTW.Range("H40").Value2 = 3833.3333
$ws_LTW.Range("I40").Value2 = 3600
$ws_LTW.Range("J40").Value2 = 5000
$ws_LTW.Range("K40").Value2 = 3600
$ws_LTW.Range("L40").Value2 = 5000
$ws_LTW.Range("M40").Value2 = -3464
$ws_LTW.Range("N40").Value2 = -5272

# LTW row 68
$ws_LTW.Range("H68").Value2 = 83335220
$ws_LTW.Range("I68").Value2 = 111112460
$ws_LTW.Range("J68").Value2 = 3466.6667
$ws_LTW.Range("K68").Value2 = 111112460
$ws_LTW.Range("L68").Value2 = 3466.6667
$ws_LTW.Range("M68").Value2 = -111111711
$ws_LTW.Range("N68").Value2 = -4964.6667

# LTW row 71
$ws_LTW.Range("H71").Value2 = 83335220
$ws_LTW.Range("I71").Value2 = 111112460
$ws_LTW.Range("J71").Value2 = 3466.6667
$ws_LTW.Range("K71").Value2 = 555562300
$ws_LTW.Range("L71").Value2 = 17333.3335
$ws_LTW.Range("M71").Value2 = -555558556
$ws_LTW.Range("N71").Value2 = -24821.3335

# LTW row 82
$ws_LTW.Range("H82").Value2 = 2389.5293
$ws_LTW.Range("I82").Value2 = 2302.1667
$ws_LTW.Range("J82").Value2 = 2437.182
$ws_LTW.Range("K82").Value2 = 2302.1667
$ws_LTW.Range("L82").Value2 = 2437.182
$ws_LTW.Range("M82").Value2 = -1941.1667
$ws_LTW.Range("N82").Value2 = -3159.182

# LTW row 85
$ws_LTW.Range("H85").Value2 = 2389.5293
$ws_LTW.Range("I85").Value2 = 2302.1667
$ws_LTW.Range("J85").Value2 = 2437.182
$ws_LTW.Range("K85").Value2 = 2302.1667
$ws_LTW.Range("L85").Value2 = 2437.182
$ws_LTW.Range("M85").Value2 = -1054.1667
$ws_LTW.Range("N85").Value2 = -4933.182

# LTW row 126
$ws_LTW.Range("H126").Value2 = 2564.9092
$ws_LTW.Range("I126").Value2 = 2629.8333
$ws_LTW.Range("J126").Value2 = 2487
$ws_LTW.Range("K126").Value2 = 7889.499899999999
$ws_LTW.Range("L126").Value2 = 7461
$ws_LTW.Range("M126").Value2 = -5419.499899999999
$ws_LTW.Range("N126").Value2 = -12401

# WVR row 122
$ws_WVR.Range("H122").Value2 = 33491.97
$ws_WVR.Range("I122").Value2 = 35560.38
$ws_WVR.Range("J122").Value2 = 3500
$ws_WVR.Range("K122").Value2 = 106681.14
$ws_WVR.Range("L122").Value2 = 10500
$ws_WVR.Range("M122").Value2 = -104231.14
$ws_WVR.Range("N122").Value2 = -15400
